$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values that look numeric are stored as text,
# matching the source data (inline strings) rather than being auto-converted
# to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.457.07"
$ws.Range("E2").Value = "  +2.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.090.13"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.88"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.666"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.19"
$ws.Range("E8").Value = "  +26.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "62.53"
$ws.Range("E9").Value = "  +1.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.384"
$ws.Range("E10").Value = "  +4.63%  "
$ws.Range("E11").Value = "  +3.52%  "
$ws.Range("E12").Value = "  +7.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.43"
$ws.Range("E13").Value = "  +5.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.389.51"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.30"
$ws.Range("E16").Value = "  +5.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.085.76"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.403.09"
$ws.Range("E18").Value = "  +2.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.44"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.78"
$ws.Range("E20").Value = "  +15.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0854"
$ws.Range("E21").Value = "  +4.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "241.37"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("E23").Value = "  +5.48%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "172.15"
$ws.Range("E26").Value = "  +1.27%  "
$ws.Range("E27").Value = "  +4.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.08"
$ws.Range("E28").Value = "  +2.44%  "
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.124"
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.24"
$ws.Range("E31").Value = "  +5.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.12"
$ws.Range("E32").Value = "  +22.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.59"
$ws.Range("E33").Value = "  +3.63%  "
$ws.Range("E34").Value = "  +6.20%  "
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("E36").Value = "  +6.93%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.31"
$ws.Range("E38").Value = "  -0.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.86"
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("E41").Value = "  +5.27%  "
$ws.Range("E42").Value = "  +10.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0992"
$ws.Range("E43").Value = "  +18.96%  "
$ws.Range("E44").Value = "  -0.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.52"
$ws.Range("E45").Value = "  +1.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.48"
$ws.Range("E46").Value = "  +118.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.80"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.330.41"
$ws.Range("E48").Value = "  -2.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.96"
$ws.Range("E49").Value = "  +4.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.39"
$ws.Range("E50").Value = "  +5.95%  "
$ws.Range("E51").Value = "  +10.41%  "
